$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.5473696196766112, 0.5473696196766112, 0.5464052287581699, 0.07889549794121917, 3.506959684504892)
    3 = @(0.5734579893152015, 0.5734579893152015, 0.5101258303180861, 0.08008931355801432, 3.796035086982494)
    4 = @(0.4995734192814485, 0.4995734192814485, 0.496, 0.1029514443045073, 3.809720685124019)
    5 = @(0.4694233531705315, 0.4694233531705315, 0.4962039045553145, 0.09396297463850045, 3.763590248074376)
    6 = @(0.5422876304023845, 0.5422876304023845, 0.5377049180327869, 0.09582944063103611, 3.484403448339702)
    7 = @(0.5754173524570891, 0.5754173524570891, 0.5993571100566355, 0.1320409098441566, 3.029877839890902)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
}
